$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-17 (columns A-J)
$data = @(
    @(0.875,              0.125,               1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.00.0", 0.66,                 0.66,                 0.66),
    @(0.875,              0.125,               1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.01.0", 0.66,                 0.66,                 0.66),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.04.0", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.05.0", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.06.0", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.07.0", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.07.1", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.07.2", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.9,                0.1,                 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.07.3", 0.6600000000000001,  0.6600000000000001,  0.6600000000000001),
    @(0.8181818181818182, 0.09090909090909091, 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.08.0", 0.6300000000000001,  0.6300000000000001,  0.6300000000000001),
    @(0.8333333333333334, 0.08333333333333333, 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.09.0", 0.6325000000000001,  0.6325000000000001,  0.6325000000000001),
    @(0.8461538461538461, 0.07692307692307693, 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.10.0", 0.6346153846153846,  0.6346153846153846,  0.6346153846153846),
    @(0.8461538461538461, 0.07692307692307693, 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.10.2", 0.6346153846153846,  0.6346153846153846,  0.6346153846153846),
    @(0.8461538461538461, 0.07692307692307693, 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.10.3", 0.6346153846153846,  0.6346153846153846,  0.6346153846153846),
    @(0.8461538461538461, 0.07692307692307693, 1,                   79.31, 8.050000000000001, "Desenvolvimento", "v1.11.0", 0.6346153846153846,  0.6346153846153846,  0.6346153846153846),
    @(0.8421052631578947, 0.05263157894736842, 0.868421052631579,   79.31, 8.050000000000001, "Desenvolvimento", "v2.00.0", 0.5818421052631579,  0.5818421052631579,  0.5818421052631579)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
